$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.636.80'
$ws.Range('E2').Value = '  +1.02%  '
$ws.Range('D3').Value = '2.247.73'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.01'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.01'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.96'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.570'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.517'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.65%  '
$ws.Range('E10').Value = '  +0.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0801'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.21'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.19%  '
$ws.Range('E13').Value = '  +0.22%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.591.41'
$ws.Range('E14').Value = '  +0.28%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.242.95'
$ws.Range('E15').Value = '  -3.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.833'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.58'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.40%  '
$ws.Range('D18').Value = '44.444.68'
$ws.Range('E18').Value = '  +1.20%  '
$ws.Range('D19').Value = '0.0₃0937'
$ws.Range('E19').Value = '  -2.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.20'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.71'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.34'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.70'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.46%  '
$ws.Range('E24').Value = '  +0.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.98'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.90%  '
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.30'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.76'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.92'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.89'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.96'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '148.22'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0784'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.28%  '
$ws.Range('E34').Value = '  +0.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.21'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.46%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.109'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.47%  '
$ws.Range('B37').Value = 'Stellar'
$ws.Range('C37').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.118'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.54%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.86'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +5.46%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '15.27'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +6.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.35'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -5.81%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.79'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.78%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('D44').Value = '1.807.32'
$ws.Range('E44').Value = '  +3.90%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.76'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +12.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '82.04'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.83%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.187'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '98.40'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.82'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '68.91'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.62%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '54.09'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.66%  '
